$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AG2").Value = "Digit_before_after"
$ws.Range("AJ2").Value = 7
$ws.Range("AG3").Value = "Digit_before_after"
$ws.Range("AJ3").Value = 6
$ws.Range("AJ4").Value = 2
$ws.Range("AG5").Value = "DICHOTIC_PRE"
$ws.Range("AJ6").Value = 8
$ws.Range("AG7").Value = "DICHOTIC_PRE"
$ws.Range("AJ7").Value = 1
$ws.Range("AG8").Value = "Digit_before_after"
$ws.Range("AJ8").Value = 5
$ws.Range("AJ9").Value = 6
$ws.Range("AG10").Value = "DICHOTIC_POST"
$ws.Range("AJ10").Value = 8
$ws.Range("AG11").Value = "Digit_before_after"
$ws.Range("AJ11").Value = 7
$ws.Range("AJ12").Value = 6
$ws.Range("AG13").Value = "DIGIT_PRE"
$ws.Range("AJ13").Value = 2
$ws.Range("AG14").Value = "DIGIT_PRE"
$ws.Range("AG15").Value = "DICHOTIC_PRE"
$ws.Range("AJ15").Value = 5
$ws.Range("AJ16").Value = 2
$ws.Range("AG17").Value = "DICHOTIC_POST"
$ws.Range("AJ17").Value = 2
$ws.Range("AG18").Value = "Dichotic_before_after"
$ws.Range("AJ18").Value = 3
$ws.Range("AG19").Value = "Dichotic_before_after"
$ws.Range("AG20").Value = "DIGIT_POST"
$ws.Range("AJ20").Value = 7
$ws.Range("AG21").Value = "Dichotic_before_after"
$ws.Range("AJ21").Value = 3
$ws.Range("AJ22").Value = 1
$ws.Range("AG23").Value = "DIGIT_PRE"
$ws.Range("AJ23").Value = 1
$ws.Range("AG24").Value = "DIGIT_POST"
$ws.Range("AJ24").Value = 1
$ws.Range("AG25").Value = "Digit_before_after"
$ws.Range("AJ25").Value = 4
$ws.Range("AG26").Value = "DIGIT_PRE"
$ws.Range("AJ26").Value = 2
$ws.Range("AG27").Value = "DICHOTIC_PRE"
$ws.Range("AJ27").Value = 1
$ws.Range("AG28").Value = "DIGIT_POST"
$ws.Range("AJ28").Value = 3
$ws.Range("AG29").Value = "Dichotic_before_after"
$ws.Range("AJ29").Value = 2
$ws.Range("AG31").Value = "DIGIT_POST"
$ws.Range("AG32").Value = "Dichotic_before_after"
$ws.Range("AJ32").Value = 4
$ws.Range("AG33").Value = "DICHOTIC_POST"
$ws.Range("AG34").Value = "DIGIT_PRE"
$ws.Range("AJ34").Value = 7
$ws.Range("AG35").Value = "Digit_before_after"
$ws.Range("AJ35").Value = 7
$ws.Range("AG36").Value = "DIGIT_PRE"
$ws.Range("AJ36").Value = 8
$ws.Range("AG37").Value = "DICHOTIC_POST"
$ws.Range("AJ37").Value = 3
$ws.Range("AG38").Value = "DIGIT_POST"
$ws.Range("AJ38").Value = 7
$ws.Range("AG39").Value = "DICHOTIC_PRE"
$ws.Range("AJ39").Value = 7
$ws.Range("AJ40").Value = 6
$ws.Range("AG41").Value = "DICHOTIC_POST"
$ws.Range("AJ41").Value = 2
$ws.Range("AG42").Value = "DIGIT_POST"
$ws.Range("AJ42").Value = 8
$ws.Range("AG43").Value = "DIGIT_POST"
$ws.Range("AJ43").Value = 2
$ws.Range("AG44").Value = "Dichotic_before_after"
$ws.Range("AJ44").Value = 1
$ws.Range("AJ45").Value = 7
$ws.Range("AG46").Value = "DICHOTIC_POST"
$ws.Range("AJ46").Value = 5
$ws.Range("AJ47").Value = 7
$ws.Range("AG48").Value = "DICHOTIC_PRE"
$ws.Range("AJ48").Value = 1
$ws.Range("AG49").Value = "Digit_before_after"
$ws.Range("AJ49").Value = 3
$ws.Range("AG50").Value = "DIGIT_PRE"
$ws.Range("AJ50").Value = 4
$ws.Range("AG51").Value = "DICHOTIC_POST"
$ws.Range("AJ51").Value = 5
$ws.Range("AG52").Value = "Digit_before_after"
$ws.Range("AG53").Value = "Dichotic_before_after"
$ws.Range("AJ53").Value = 3
$ws.Range("AG54").Value = "DICHOTIC_PRE"
$ws.Range("AJ54").Value = 2
$ws.Range("AG55").Value = "DIGIT_POST"
$ws.Range("AJ55").Value = 5
$ws.Range("AG56").Value = "Dichotic_before_after"
$ws.Range("AJ56").Value = 1
$ws.Range("AG57").Value = "Digit_before_after"
$ws.Range("AJ57").Value = 6
$ws.Range("AG58").Value = "DICHOTIC_POST"
$ws.Range("AJ58").Value = 1
$ws.Range("AG59").Value = "DICHOTIC_PRE"
$ws.Range("AJ59").Value = 8
$ws.Range("AG60").Value = "Digit_before_after"
$ws.Range("AJ60").Value = 4
$ws.Range("AG61").Value = "DIGIT_PRE"
$ws.Range("AJ61").Value = 1
$ws.Range("AG62").Value = "DIGIT_PRE"
$ws.Range("AJ62").Value = 6
$ws.Range("AG63").Value = "DIGIT_PRE"
$ws.Range("AJ63").Value = 5
$ws.Range("AJ64").Value = 8
$ws.Range("AG65").Value = "Dichotic_before_after"
$ws.Range("AG66").Value = "DIGIT_POST"
$ws.Range("AJ66").Value = 1
$ws.Range("AG67").Value = "DIGIT_PRE"
$ws.Range("AJ67").Value = 4
$ws.Range("AG68").Value = "DICHOTIC_PRE"
$ws.Range("AJ68").Value = 2
$ws.Range("AG69").Value = "DICHOTIC_PRE"
$ws.Range("AJ69").Value = 4
$ws.Range("AG70").Value = "DICHOTIC_POST"
$ws.Range("AJ70").Value = 5
$ws.Range("AG71").Value = "DIGIT_PRE"
$ws.Range("AJ71").Value = 3
$ws.Range("AG72").Value = "Dichotic_before_after"
$ws.Range("AJ72").Value = 6
$ws.Range("AG73").Value = "DICHOTIC_PRE"
$ws.Range("AG74").Value = "DICHOTIC_POST"
$ws.Range("AJ74").Value = 2
$ws.Range("AG75").Value = "DIGIT_POST"
$ws.Range("AJ75").Value = 7
$ws.Range("AG76").Value = "Dichotic_before_after"
$ws.Range("AJ76").Value = 7
$ws.Range("AJ77").Value = 4
$ws.Range("AG78").Value = "Digit_before_after"
$ws.Range("AJ78").Value = 5
$ws.Range("AG79").Value = "DICHOTIC_PRE"
$ws.Range("AJ79").Value = 4
$ws.Range("AJ80").Value = 1
$ws.Range("AG81").Value = "DIGIT_PRE"
$ws.Range("AJ81").Value = 7
